$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.023.12"
$ws.Range("E2").Value = "  -7.06%  "
$ws.Range("D3").Value = "3.259.09"
$ws.Range("E3").Value = "  -8.56%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "549.80"
$ws.Range("E5").Value = "  -7.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.93"
$ws.Range("E6").Value = "  -10.40%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.583"
$ws.Range("E8").Value = "  -5.08%  "
$ws.Range("D9").Value = "3.252.35"
$ws.Range("E9").Value = "  -8.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.183"
$ws.Range("E10").Value = "  -12.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.578"
$ws.Range("E11").Value = "  -8.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.60"
$ws.Range("E12").Value = "  -12.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000259"
$ws.Range("E13").Value = "  -10.84%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.43"
$ws.Range("E14").Value = "  -9.57%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "3.780.45"
$ws.Range("E15").Value = "  -8.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "601.29"
$ws.Range("E16").Value = "  -8.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "17.75"
$ws.Range("E17").Value = "  -4.21%  "
$ws.Range("D18").Value = "64.849.18"
$ws.Range("E18").Value = "  -7.39%  "
$ws.Range("E19").Value = "  -4.51%  "
$ws.Range("D20").Value = "3.267.73"
$ws.Range("E20").Value = "  -8.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.27"
$ws.Range("E21").Value = "  -11.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.891"
$ws.Range("E22").Value = "  -8.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.36"
$ws.Range("E23").Value = "  -5.64%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "101.37"
$ws.Range("E24").Value = "  -3.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.92"
$ws.Range("E25").Value = "  -9.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.92"
$ws.Range("E26").Value = "  -11.47%  "
$ws.Range("E27").Value = "  -1.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.64"
$ws.Range("E28").Value = "  -10.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.20"
$ws.Range("E29").Value = "  -10.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.54"
$ws.Range("E30").Value = "  -11.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "30.07"
$ws.Range("E31").Value = "  -10.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.79"
$ws.Range("E32").Value = "  -13.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.15"
$ws.Range("E33").Value = "  -10.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.90"
$ws.Range("E34").Value = "  -7.91%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.103"
$ws.Range("E35").Value = "  -7.73%  "
$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D36").Value = "3.731.72"
$ws.Range("E36").Value = "  -0.67%  "
$ws.Range("E37").Value = "  -0.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "526.89"
$ws.Range("E38").Value = "  +2.43%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "55.66"
$ws.Range("E39").Value = "  -10.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.40"
$ws.Range("E40").Value = "  -7.71%  "
$ws.Range("D41").Value = "0.0₃0696"
$ws.Range("E41").Value = "  -14.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.63"
$ws.Range("E42").Value = "  -11.86%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.124"
$ws.Range("E43").Value = "  -8.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "31.45"
$ws.Range("E44").Value = "  -10.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.333"
$ws.Range("E45").Value = "  -11.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.17"
$ws.Range("E46").Value = "  +16.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.17"
$ws.Range("E47").Value = "  -7.27%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0403"
$ws.Range("E48").Value = "  -11.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.128"
$ws.Range("E49").Value = "  -7.09%  "
$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.55"
$ws.Range("E50").Value = "  -11.97%  "
$ws.Range("B51").Value = "FirstDigitalUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.997"
$ws.Range("E51").Value = "  -0.33%  "
